$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Step 1: the old row 16 (user/addGallery) had F16 = "Inactive".
# The author repurposed that status text to "Active" before adding the
# new "user/active" service row above it.
$ws.Cells.Item(16, 6).Value = "Active"

# --- Step 2: insert a brand-new row above the old row 16, shifting
# everything (old rows 16-34) down by one (to 17-35).
$ws.Rows.Item(16).Insert()

# --- Step 3: copy the formatting (styles) of row 12 ("user/register"),
# which uses the alternate banding color, into the new row 16.
$ws.Range("A12:F12").Copy()
$ws.Range("A16:F16").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Rows.Item(16).RowHeight = 30

# --- Step 4: populate the new "user/active" service row.
$ws.Cells.Item(16, 1).Value = "user/active"
$ws.Cells.Item(16, 2).Value = "{}"
$ws.Cells.Item(16, 3).Value = "JSON"
$ws.Cells.Item(16, 4).Value = '{ "data" : ["id" : id/-1] , "error" : [] }'
$ws.Cells.Item(16, 5).Value = $ws.Cells.Item(14, 5).Value2
$ws.Cells.Item(16, 6).Value = "Active"

# --- Step 5: extend the conditional formatting ("Active"/"Inactive"
# coloring) so it still covers the whole status column after the insert.
$cf = $ws.Range("F11:F34").FormatConditions
$cf.Item(1).ModifyAppliesToRange($ws.Range("F11:F35"))

# --- Step 6: refresh the view selection to match the new layout.
$excel.ActiveWindow.ScrollRow = 11
$ws.Range("F18").Select()
